$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(46048, 12095.1033375746, 11618.8480441374, 19995.86, 7967.53177860242, -17.0616740525061),
    @(46049, 12216.8447580457, 11995.5718674379, 12075.86, 8222.76217015504, 339.269751566374),
    @(46050, 12156.7645803242, 11524.1578089809, 12075.86, 8164.98106466119, 317.219953068422),
    @(46051, 12068.8906703456, 11603.1782668546, 12075.86, 8100.21288022034, 317.813797794788),
    @(46052, 11952.3080730186, 10896.3462944393, 12075.86, 8039.62913860075, 285.838143043335),
    @(46053, 4866.94803138697, 7810.59150886682, 12075.86, 7929.96897791077, 152.695853615733),
    @(46054, 5046.18659312722, 7836.21906380889, 9743.86, 7892.57678860208, 249.37232718379),
    @(46055, 11445.1451222577, 11171.4994288684, 9743.86, 7808.42406784652, 384.835979029789),
    @(46056, 11445.1451222577, 10751.8757132856, 9743.86, 7808.42406784652, 367.351657547171),
    @(46057, 11445.1451222577, 10402.3851065668, 9743.86, 7808.42406784652, 352.789548933889),
    @(46058, 11445.1451222577, 10493.7628986866, 9743.86, 7808.42406784652, 356.596956938882),
    @(46059, 11445.1451222577, 9815.93594729526, 9743.86, 7808.42406784652, 328.354167297574),
    @(46060, 4816.06445810594, 6974.19612391431, 9743.86, 7435.45854081644, 194.408111030448),
    @(46061, 4716.67601641354, 7029.01798450083, 9743.86, 7426.7951422196, 196.331380280018)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
